# Generate Report for Handback
# Marks the two still-outstanding files (98cda078... and e9225a4f...) as
# handed back (in sync with en-US) for both the zh-cn and de-de locales,
# filling in their "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns, and mirrors the status onto the
# Overview sheet.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: rows 4 & 5 (98cda078... / e9225a4f...) move from
# "Ready for handoff" to "Handed back: in sync with en-US" for both the
# zh-cn (E) and de-de (F) columns.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = $handedBack
$wsOverview.Range("F4").Value = $handedBack
$wsOverview.Range("E5").Value = $handedBack
$wsOverview.Range("F5").Value = $handedBack

# ---------------------------------------------------------------------
# zh-cn sheet: update Status, Latest Target File, Latest Handback File
# and Latest Handback DateTime for rows 4 & 5.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C4").Value = $handedBack
$g4 = $wsZhCn.Range("G4").Value2
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0136ab465805189892dbff124a7e0a5e0a3c509/e2e/98cda078-543e-48d5-b448-e610ef084672.md", [Type]::Missing, [Type]::Missing, "98cda078-543e-48d5-b448-e610ef084672.md") | Out-Null
$wsZhCn.Range("J4").Value = $g4
$wsZhCn.Range("K4").Value = "2016-08-31 20:32:56"

$wsZhCn.Range("C5").Value = $handedBack
$g5 = $wsZhCn.Range("G5").Value2
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0136ab465805189892dbff124a7e0a5e0a3c509/e2e/e9225a4f-dd7f-45e7-8d63-b79168e467fc.md", [Type]::Missing, [Type]::Missing, "e9225a4f-dd7f-45e7-8d63-b79168e467fc.md") | Out-Null
$wsZhCn.Range("J5").Value = $g5
$wsZhCn.Range("K5").Value = "2016-08-31 20:32:56"

# ---------------------------------------------------------------------
# de-de sheet: same updates, with its own handback timestamp.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C4").Value = $handedBack
$g4d = $wsDeDe.Range("G4").Value2
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0136ab465805189892dbff124a7e0a5e0a3c509/e2e/98cda078-543e-48d5-b448-e610ef084672.md", [Type]::Missing, [Type]::Missing, "98cda078-543e-48d5-b448-e610ef084672.md") | Out-Null
$wsDeDe.Range("J4").Value = $g4d
$wsDeDe.Range("K4").Value = "2016-08-31 20:33:12"

$wsDeDe.Range("C5").Value = $handedBack
$g5d = $wsDeDe.Range("G5").Value2
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0136ab465805189892dbff124a7e0a5e0a3c509/e2e/e9225a4f-dd7f-45e7-8d63-b79168e467fc.md", [Type]::Missing, [Type]::Missing, "e9225a4f-dd7f-45e7-8d63-b79168e467fc.md") | Out-Null
$wsDeDe.Range("J5").Value = $g5d
$wsDeDe.Range("K5").Value = "2016-08-31 20:33:12"

Write-Host "Handback report generated."
